$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("A1").Value = "TEST"
Write-Output "done"
